$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for "Haba" that belongs
# chronologically between the existing rows 66 and 67 (old numbering).
# Insert a fresh row at position 67, which pushes the previous rows
# 67..105 down to 68..106 (matches the dimension growing from
# A1:R105 to A1:R106 in the diff).
$ws.Rows.Item(67).Insert()

# Populate the newly inserted row 67 with the new record's data.
$ws.Range("A67").Value = 6
$ws.Range("B67").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C67").Value = "Metropolitana"
$ws.Range("D67").Value = 44438
$ws.Range("E67").Value = 13
$ws.Range("F67").Value = 100112026
$ws.Range("G67").Value = "Haba"
$ws.Range("H67").Value = "Sin especificar"
$ws.Range("I67").Value = "Primera"
$ws.Range("J67").Value = 400
$ws.Range("K67").Value = 15000
$ws.Range("L67").Value = 16000
$ws.Range("M67").Value = 15575
$ws.Range("N67").Value = "`$/saco 25 kilos"
$ws.Range("O67").Value = "Provincia de Huasco"
$ws.Range("P67").Value = 623
$ws.Range("Q67").Value = 25
$ws.Range("R67").Value = "Hortaliza"
